$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '80.230.18'
$ws.Range('E2').Value = '  +5.33%  '
$ws.Range('D3').Value = '3.213.95'
$ws.Range('E3').Value = '  +4.80%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.274'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +32.10%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.88%  '
$ws.Range('D10').Value = '3.215.13'
$ws.Range('E10').Value = '  +4.84%  '
$ws.Range('E11').Value = '  +42.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000260'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +34.89%  '
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').Value = '3.802.86'
$ws.Range('E15').Value = '  +5.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +11.99%  '
$ws.Range('D17').Value = '80.272.13'
$ws.Range('E17').Value = '  +5.41%  '
$ws.Range('D18').Value = '3.216.89'
$ws.Range('E18').Value = '  +4.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +17.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +20.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +22.59%  '
$ws.Range('D24').Value = '3.382.45'
$ws.Range('E24').Value = '  +4.79%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '77.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.40%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +12.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.997'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000124'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +15.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.58%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.29%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '555.56'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +11.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.153'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +24.03%  '
$ws.Range('E35').Value = '  +6.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +14.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.125'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +22.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.421'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '20.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.76'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +13.11%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '191.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +12.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +13.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.796'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.34'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.640'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.88%  '
